$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 29, shifting existing rows 29-32 down to 30-33
$ws.Rows.Item(29).Insert()

# Fill in the new row 29 with the new record's data
$ws.Cells.Item(29, 1).Value = 10
$ws.Cells.Item(29, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(29, 3).Value = "La Araucanía"
$ws.Cells.Item(29, 4).Value = 44511
$ws.Cells.Item(29, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(29, 5).Value = 9
$ws.Cells.Item(29, 6).Value = 100114002
$ws.Cells.Item(29, 7).Value = "Camote"
$ws.Cells.Item(29, 8).Value = "Sin especificar"
$ws.Cells.Item(29, 9).Value = "Primera"
$ws.Cells.Item(29, 10).Value = 50
$ws.Cells.Item(29, 11).Value = 20000
$ws.Cells.Item(29, 12).Value = 20000
$ws.Cells.Item(29, 13).Value = 20000
$ws.Cells.Item(29, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(29, 15).Value = "Perú"
$ws.Cells.Item(29, 16).Value = 1000
$ws.Cells.Item(29, 17).Value = 20
$ws.Cells.Item(29, 18).Value = "Hortaliza"
